$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: I1 = "I0", J1 = "IF" (same style as the other header cells,
# e.g. H1). Copy the formatting from H1 so the new cells reuse the exact
# same cell style rather than creating a brand-new one.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# Data rows 2-17: I column and J column values
$values = @{
    2  = @(8, 8)
    3  = @(8, 8)
    4  = @(9, 9)
    5  = @(8, 8)
    6  = @(9, 9)
    7  = @(9, 9)
    8  = @(6, 7)
    9  = @(9, 9)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(8, 8)
    13 = @(9, 9)
    14 = @(8, 8)
    15 = @(7, 7)
    16 = @(7, 7)
    17 = @(7, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
